$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (re-check results)
$ws.Range("B2").Value = "2024-06-16 10:50:22.225205"
$ws.Range("C2").Value = 21399

# Update existing row 3 (re-check results)
$ws.Range("B3").Value = "2024-06-16 10:50:46.586880"
$ws.Range("C3").Value = 21399
$ws.Range("D3").Value = 0

# New row 4 -- copy the formatting used by column A's index cells (style id 1)
# so the new index cell matches the existing bold/bordered/centered look.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2024-06-16 12:49:49.275805"
$ws.Range("C4").Value = 21402
$ws.Range("D4").Value = 3

# New row 5
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2024-06-16 12:50:15.873938"
$ws.Range("C5").Value = 21404
$ws.Range("D5").Value = 2

$excel.CutCopyMode = 0
